$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The schema table shrank from 21 data rows to 16 (fields were dropped,
# merged, reordered and a few new ones added) -- drop the trailing rows first.
$ws.Rows("17:21").Delete()

# Rewrite the Field / Description / Type table in place to match the new schema.
$ws.Range("A1").Value = 'Field'
$ws.Range("B1").Value = 'Description'
$ws.Range("C1").Value = 'Type'
$ws.Range("A2").Value = 'dataset'
$ws.Range("B2").Value = 'GEO series id(s), if provided, in python list format. e.g. GSE190147'
$ws.Range("C2").Value = 'list'
$ws.Range("A3").Value = 'pmid'
$ws.Range("B3").Value = 'PubMed ID(s), if provided, in python list format. e.g. 35926038'
$ws.Range("C3").Value = 'list'
$ws.Range("A4").Value = 'pmcid'
$ws.Range("B4").Value = 'PMC ID(s), if provided, in python list format. e.g. PMC9371140'
$ws.Range("C4").Value = 'list'
$ws.Range("A5").Value = 'doi'
$ws.Range("B5").Value = 'Publication doi(s), without https://head, if provided, in python list format. e.g. 10.1126/science.abn5800'
$ws.Range("C5").Value = 'list'
$ws.Range("A6").Value = 'other_ids'
$ws.Range("B6").Value = 'INSDC or other project IDs, in python list format.'
$ws.Range("C6").Value = 'list'
$ws.Range("A7").Value = 'title'
$ws.Range("B7").Value = 'Project title, if provided.'
$ws.Range("C7").Value = 'str'
$ws.Range("A8").Value = 'project_description'
$ws.Range("B8").Value = 'Description or abstrcut of what this project studied.'
$ws.Range("C8").Value = 'str'
$ws.Range("A9").Value = 'species'
$ws.Range("B9").Value = 'Specie(s) involved in this study, one or multiple, denoted by Latin scientific name, in python list format.'
$ws.Range("C9").Value = 'list'
$ws.Range("A10").Value = 'organ'
$ws.Range("B10").Value = 'Organisation(s) of sample sampling in the study, e.g. `boold`, `lung`. Store as python list format.'
$ws.Range("C10").Value = 'list'
$ws.Range("A11").Value = 'topic'
$ws.Range("B11").Value = 'Topic of given paper, select one or multiple from `aging`, `bioinformaticstool`, `brainorganoids`, `cancer`, `cardiology`, `cellbiology`, `developmentalbiology`, `disease`, `drugaddiction`, `epigenetics`, `fibrosis`, `hematopoiesis`, `host-viralinteractions`, `immunology`, `infectiousdisease`, `inflammation`, `metabolism`, `methodspaper`, `multi-omics`, `neuroscience`, `organtransplantation`, `pathogenesis`, `proteomics`, `psychiatricdisorders`, `pulmonology`, `regeneration`, `relhomeostasis`, `review`, `sensory`, `single-cellanalysis`, `stemcells`, `Tcellbiology`, `tissueengineering`, `transcriptomics`, or `Noneofabove`. Store as python list format.'
$ws.Range("C11").Value = 'list'
$ws.Range("A12").Value = 'resolution'
$ws.Range("B12").Value = 'The study resolution of this research, `Single-cell` for single-cell/single-nucleus research using technics such as 10x genomics 3'' RNA-seq; `Spatial` for special omics tech such as 10x visum; `Bulk` for bulk level research; or `NS` if not speicific.'
$ws.Range("C12").Value = 'str'
$ws.Range("A13").Value = 'technology_name'
$ws.Range("B13").Value = 'Name of single cell sequenceing technology used in project, e.g. `SMART-seq2`, `10x genomics chromium single cell 3''` (or `10x 3'' V3` for short, depent on the authors description), or `NS` if not specific. In python list format. Please extract the author''s original statements.'
$ws.Range("C13").Value = 'list'
$ws.Range("A14").Value = 'disease'
$ws.Range("B14").Value = 'What kind of disease involved in this research, in python list format, and use `Normal` for indicate non-diease samples were involved.'
$ws.Range("C14").Value = 'list'
$ws.Range("A15").Value = 'library_strategy'
$ws.Range("B15").Value = 'Seqencing library strategy, if given. e.g. `RNA-Seq`, `ATAC-Seq`, `DNA-Seq` etc, in python list format.'
$ws.Range("C15").Value = 'list'
$ws.Range("A16").Value = 'nuclei_extraction'
$ws.Range("B16").Value = 'Whether it is using nuclei only for single-cell level resolution sequencing. TRUE`, `FALSE` or `NS` if no hint. e.g. snRNA-seq is `TRUE`.'
$ws.Range("C16").Value = 'str'

# Column B (Description) is now wider to fit the longer text.
$ws.Columns("B").ColumnWidth = 32.29

# Scroll position / selection as left by the author on save.
$ws.Range("H17").Select()
